$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and B13/B14, C13/C14 swap)
$updates = @(
    @{Cell='D2'; Value='30.035.55'},
    @{Cell='E2'; Value='  +0.88%  '},
    @{Cell='D3'; Value='1.885.58'},
    @{Cell='E3'; Value='  +1.25%  '},
    @{Cell='D4'; Value='0.9986'},
    @{Cell='E4'; Value='  -0.56%  '},
    @{Cell='D5'; Value='0.7368'},
    @{Cell='E5'; Value='  +0.68%  '},
    @{Cell='D6'; Value='242.10'},
    @{Cell='E6'; Value='  +0.35%  '},
    @{Cell='D7'; Value='0.9992'},
    @{Cell='D8'; Value='0.3166'},
    @{Cell='E8'; Value='  +2.77%  '},
    @{Cell='D9'; Value='0.07188'},
    @{Cell='E9'; Value='  +2.37%  '},
    @{Cell='D10'; Value='24.77'},
    @{Cell='E10'; Value='  +1.55%  '},
    @{Cell='D11'; Value='0.08331'},
    @{Cell='E11'; Value='  -0.94%  '},
    @{Cell='D12'; Value='0.7571'},
    @{Cell='E12'; Value='  +1.67%  '},
    @{Cell='B13'; Value='WrappedEther'},
    @{Cell='C13'; Value='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'},
    @{Cell='D13'; Value='1.909.67'},
    @{Cell='E13'; Value='  +1.97%  '},
    @{Cell='B14'; Value='Polkadot'},
    @{Cell='C14'; Value='https://coinranking.com/coin/25W7FG7om+polkadot-dot'},
    @{Cell='D14'; Value='5.403'},
    @{Cell='E14'; Value='  +1.95%  '},
    @{Cell='D15'; Value='93.01'},
    @{Cell='E15'; Value='  +1.21%  '},
    @{Cell='D16'; Value='6.161'},
    @{Cell='E16'; Value='  +1.70%  '},
    @{Cell='D17'; Value='30.037.10'},
    @{Cell='E17'; Value='  +0.82%  '},
    @{Cell='D18'; Value='249.93'},
    @{Cell='E18'; Value='  +4.56%  '},
    @{Cell='D19'; Value='13.59'},
    @{Cell='E19'; Value='  +1.07%  '},
    @{Cell='D20'; Value='0.000007857'},
    @{Cell='E20'; Value='  +1.44%  '},
    @{Cell='D21'; Value='2.144.20'},
    @{Cell='E21'; Value='  -0.20%  '},
    @{Cell='D22'; Value='0.9986'},
    @{Cell='E22'; Value='  -0.34%  '},
    @{Cell='D23'; Value='7.894'},
    @{Cell='E23'; Value='  +0.19%  '},
    @{Cell='D24'; Value='0.9983'},
    @{Cell='E24'; Value='  -0.58%  '},
    @{Cell='D25'; Value='0.1562'},
    @{Cell='E25'; Value='  +0.16%  '},
    @{Cell='D26'; Value='9.281'},
    @{Cell='E26'; Value='  +0.47%  '},
    @{Cell='D27'; Value='163.27'},
    @{Cell='E27'; Value='  +0.81%  '},
    @{Cell='D28'; Value='18.68'},
    @{Cell='E28'; Value='  +1.24%  '},
    @{Cell='D29'; Value='2.050'},
    @{Cell='E29'; Value='  +2.97%  '},
    @{Cell='E30'; Value='  -0.69%  '},
    @{Cell='D31'; Value='4.570'},
    @{Cell='E31'; Value='  +3.26%  '},
    @{Cell='E32'; Value='  +0.34%  '},
    @{Cell='D33'; Value='4.201'},
    @{Cell='E33'; Value='  +2.14%  '},
    @{Cell='D34'; Value='0.05335'},
    @{Cell='E34'; Value='  -0.04%  '},
    @{Cell='D35'; Value='1.250'},
    @{Cell='E35'; Value='  +2.16%  '},
    @{Cell='D36'; Value='0.7693'},
    @{Cell='E36'; Value='  +3.91%  '},
    @{Cell='D37'; Value='0.9982'},
    @{Cell='E37'; Value='  -0.36%  '},
    @{Cell='D38'; Value='2.719'},
    @{Cell='E38'; Value='  +0.73%  '},
    @{Cell='D39'; Value='0.01962'},
    @{Cell='E39'; Value='  +2.31%  '},
    @{Cell='D40'; Value='2.757'},
    @{Cell='E40'; Value='  +0.67%  '},
    @{Cell='D41'; Value='0.4582'},
    @{Cell='E41'; Value='  +4.02%  '},
    @{Cell='E42'; Value='  +0.69%  '},
    @{Cell='D43'; Value='1.086.50'},
    @{Cell='E43'; Value='  -0.86%  '},
    @{Cell='D44'; Value='72.34'},
    @{Cell='E44'; Value='  +1.27%  '},
    @{Cell='D45'; Value='0.8729'},
    @{Cell='E45'; Value='  +1.58%  '},
    @{Cell='E46'; Value='  +2.84%  '},
    @{Cell='E47'; Value='  -0.39%  '},
    @{Cell='D48'; Value='1.857'},
    @{Cell='E48'; Value='  +2.03%  '},
    @{Cell='D49'; Value='7.585'},
    @{Cell='E49'; Value='  -0.75%  '},
    @{Cell='D50'; Value='9.563'},
    @{Cell='E50'; Value='  -0.82%  '},
    @{Cell='D51'; Value='2.047.47'},
    @{Cell='E51'; Value='  +0.44%  '}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell -match "^[DE]") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}

Write-Output ("Applied " + $updates.Count + " cell updates.")
